# Fix the description of how the RPI receives power (table cell T4.2):
# old: "De RPI wordt gevoed vanuit de Raspberry PI officiële voeding voor de 5V 5A
#       volledig vermogen modus, of vanuit de printplaat naar de GPIO pinnen gevoed
#       te worden op 5V ?A max."
# new: "De RPI word vanuit de printplaat gevoed via de GPIO pinnen."

$d = $word.ActiveDocument

$old = "De RPI wordt gevoed vanuit de Raspberry PI officiële voeding voor de 5V 5A volledig vermogen modus, of vanuit de printplaat naar de GPIO pinnen gevoed te worden op 5V ?A max."
$new = "De RPI word vanuit de printplaat gevoed via de GPIO pinnen."

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Target sentence not found - cannot apply edit."
}

Write-Output "Replaced sentence. Found=$found"
